# fix: adding total stock value
#
# Adds a new "ValorTotal" header to the EstoqueAtual sheet (column F),
# reusing the same header formatting (bold / bordered / centered) as the
# rest of row 1, and restores the last-recorded cell selections on both
# sheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("EstoqueAtual")
$ws2 = $wb.Worksheets.Item("HistoricoTransacoes")

# New header cell, formatted like its neighbour (E1: ValorUnitario/NomeProduto
# style header row), then filled with the new column's title.
$ws1.Range("E1").Copy()
$ws1.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("F1").Value = "ValorTotal"

# Restore the selections left behind on each sheet.
$ws1.Rows("2").Select()
$ws2.Range("J22").Select()
